# Added ifoCAST full series evaluation: shift each row's error series one
# quarter to the left (dropping the oldest Q0 value) and, where new data is
# available (rows 2-6), append a freshly evaluated value in the last column.
# Rows further down the staircase simply lose their last (now out-of-range)
# value with no replacement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

$rowData = @{
    2  = @(-0.4989491018657047, -0.6407780486049788, 0.002196572797085183, 0.1266092791266676, -0.2710209269109496, -0.3471896406117954, -0.07086758248887381, -0.2420874048005978, -0.08893999929185659, -0.3663687737149753)
    3  = @(-0.05818351157133772, 0.5847911098307261, 0.7092038161603086, 0.3115736101226914, 0.2354048964218456, 0.5117269545447671, 0.3405071322330432, 0.4936545377417844, 0.2162257633186657, 0.05323806995971928)
    4  = @(0.4891160690213684, 0.6135287753509509, 0.2158985693133336, 0.1397298556124878, 0.4160519137354094, 0.2448320914236854, 0.3979794969324266, 0.1205507225093079, -0.04243697084963852, -0.1665195462441563)
    5  = @(0.8739478123496736, 0.4763176063120564, 0.4001488926112106, 0.6764709507341322, 0.5052511284224082, 0.6583985339311494, 0.3809697595080307, 0.2179820661490843, 0.09389949075456649, 0.6100525277605273)
    6  = @(-0.3591373898464348, -0.4353061035472806, -0.158984045424359, -0.330203867736083, -0.1770564622273418, -0.4544852366504605, -0.617472930009407, -0.7415555054039247, -0.2254024683979639, -0.4353061035472806)
    7  = @(-0.3352267436446591, -0.0589046855217375, -0.2301245078334615, -0.07697710232472027, -0.354405876747839, -0.5173935701067854, -0.6414761455013032, -0.1253231084953424, -0.3352267436446591)
    8  = @(0.2948818205579588, 0.1236619982462347, 0.276809403754976, -0.0006193706681427817, -0.1636070640270892, -0.287689639421607, 0.2284633975843539, 0.01855976243503714)
    9  = @(-0.0245153921862106, 0.1286320133225306, -0.1487967611005881, -0.3117844544595345, -0.4358670298540523, 0.08028600715190851, -0.1296176279974082)
    10 = @(-0.02881397568162436, -0.3062427501047431, -0.4692304434636895, -0.5933130188582073, -0.07715998185224648, -0.2870636170015632)
    11 = @(0.1944167064213277, 0.0314290130623813, -0.09265356233213651, 0.4234994746738243, 0.2135958395245076)
    12 = @(-0.2489288674730878, -0.3730114428676057, 0.1431415941383551, -0.06676204101096155)
    13 = @(-0.201036585022594, 0.3151164519833668, 0.1052128168340501)
    14 = @(0.009253912237035311, -0.2006497229122814)
    15 = @(0.4116802297750048)
    16 = @()
}

foreach ($r in 2..16) {
    $values = $rowData[$r]

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $cellRef = "$colLetter$r"
        if ($i -lt $values.Length) {
            $ws.Range($cellRef).Value = $values[$i]
        } else {
            $ws.Range($cellRef).ClearContents()
        }
    }
}

$wb.Save()
